$d = $word.ActiveDocument

function Replace-Text($find, $replace) {
    $d.Content.Find.Execute($find, $true, $false, $false, $false, $false, $true, 1, $false, $replace, 2)
}

Replace-Text "2024-10-26 Saturday" "2024-10-27 Sunday"

Replace-Text "43×60=" "15×56="
Replace-Text "88×42=" "74×37="
Replace-Text "41×87=" "69×21="
Replace-Text "25×65=" "19×17="
Replace-Text "28×27=" "30×78="

Replace-Text "38×93=" "31×91="
Replace-Text "31×53=" "40×56="
Replace-Text "39×19=" "64×53="
Replace-Text "15×93=" "76×12="
Replace-Text "51×61=" "50×41="

Replace-Text "26×93=" "22×32="
Replace-Text "93×36=" "91×62="
Replace-Text "13×74=" "18×25="
Replace-Text "65×30=" "67×37="
Replace-Text "58×85=" "22×55="

Replace-Text "89×81=" "70×56="
Replace-Text "63×13=" "98×58="
Replace-Text "83×52=" "19×61="
Replace-Text "47×74=" "24×90="
Replace-Text "15×50=" "78×76="

Replace-Text "20×81=" "66×71="
Replace-Text "92×73=" "41×56="
Replace-Text "59×14=" "63×64="
Replace-Text "30×47=" "44×73="
Replace-Text "42×35=" "67×81="
